$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header for new column G
$ws.Range("G1").Value = "Time"

# Column G width
$ws.Columns.Item(7).ColumnWidth = 12.5546875

# Dates for G2:G11 (first day of each month, Jan 2025 - Oct 2025)
$dates = @("2025-01-01","2025-02-01","2025-03-01","2025-04-01","2025-05-01","2025-06-01","2025-07-01","2025-08-01","2025-09-01","2025-10-01")
for ($i = 0; $i -lt $dates.Length; $i++) {
    $row = $i + 2
    $cell = $ws.Cells.Item($row, 7)
    $cell.Value = $dates[$i]
    $cell.NumberFormat = "m/d/yyyy"
}

$ws.Range("J23").Select()
